$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'66.730.93"
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +0.15%  '
$ws.Range('D3').Value = "'3.232.89"
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +0.98%  '
$ws.Range('D4').Value = "'0.999"
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('D5').Value = "'608.79"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.77%  '
$ws.Range('D6').Value = "'158.76"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +2.44%  '
$ws.Range('E7').Value = '  +0.07%  '
$ws.Range('D8').Value = "'3.231.29"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +0.95%  '
$ws.Range('D9').Value = "'0.551"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +1.27%  '
$ws.Range('E10').Value = '  +0.38%  '
$ws.Range('D11').Value = "'5.74"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -4.65%  '
$ws.Range('D12').Value = "'0.505"
$ws.Range('D12').Style = 'Normal'
$ws.Range('D13').Value = "'0.0000272"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +1.43%  '
$ws.Range('D14').Value = "'38.90"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -0.56%  '
$ws.Range('D15').Value = "'3.764.46"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +0.99%  '
$ws.Range('D16').Value = "'66.732.93"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +0.29%  '
$ws.Range('D17').Value = "'7.38"
$ws.Range('D17').Style = 'Normal'
$ws.Range('D18').Value = "'3.237.20"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +1.10%  '
$ws.Range('E19').Value = '  +1.32%  '
$ws.Range('D20').Value = "'511.10"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -0.45%  '
$ws.Range('D21').Value = "'15.23"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.93%  '
$ws.Range('D22').Value = "'0.735"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.91%  '
$ws.Range('D23').Value = "'8.02"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.99%  '
$ws.Range('D24').Value = "'14.64"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -3.11%  '
$ws.Range('D25').Value = "'85.09"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -0.70%  '
$ws.Range('E26').Value = '  +0.14%  '
$ws.Range('E27').Value = '  -0.26%  '
$ws.Range('D28').Value = "'9.12"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -2.58%  '
$ws.Range('D29').Value = "'2.36"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +2.10%  '
$ws.Range('D30').Value = "'0.130"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +45.03%  '
$ws.Range('D31').Value = "'2.95"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +0.39%  '
$ws.Range('D32').Value = "'6.99"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -2.44%  '
$ws.Range('D33').Value = "'28.19"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -0.49%  '
$ws.Range('E34').Value = '  +0.13%  '
$ws.Range('D35').Value = "'1.18"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -4.25%  '
$ws.Range('D36').Value = "'6.51"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -0.54%  '
$ws.Range('D37').Value = "'502.99"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -0.52%  '
$ws.Range('D38').Value = "'55.57"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +1.18%  '
$ws.Range('D39').Value = "'0.0₃0769"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +11.67%  '
$ws.Range('B40').Value = 'VeChain'
$ws.Range('C40').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D40').Value = "'0.0423"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -0.15%  '
$ws.Range('B41').Value = 'dogwifhat'
$ws.Range('C41').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D41').Value = "'3.06"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +6.65%  '
$ws.Range('D42').Value = "'0.130"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +5.61%  '
$ws.Range('D43').Value = "'8.74"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -1.60%  '
$ws.Range('D44').Value = "'0.299"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.78%  '
$ws.Range('D45').Value = "'2.45"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -0.60%  '
$ws.Range('D46').Value = "'2.913.16"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -0.71%  '
$ws.Range('D47').Value = "'28.16"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -1.47%  '
$ws.Range('D48').Value = "'2.42"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +3.86%  '
$ws.Range('E49').Value = '  -0.04%  '
$ws.Range('D50').Value = "'0.116"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -0.70%  '
$ws.Range('D51').Value = "'122.83"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +0.17%  '
